# CIERRE 6 OCT 22
# - Update the incentive-month label from AGOSTO to SEPTIEMBRE on the
#   "VALES DE INSENTIVOS" sheet.
# - Make "VALES DE INSENTIVOS" the active/selected tab (it was
#   "ARQUITECTO" before).
# (The TODAY() cells recalc automatically against the replay clock, so
# they are left alone.)

$wb = $excel.ActiveWorkbook

$wsVales = $wb.Worksheets.Item("VALES DE INSENTIVOS")
$wsVales.Range("A4").Value = "PAGO DE INCENTIVO DEL MES DE  SEPTIEMBRE   2022"

$wsVales.Activate()
